$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.809.37"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.805.38"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "4.446.74"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "3.830.33"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "67.833.10"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D29").Value = "3.956.89"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.89%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "391.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "
